$wb = $excel.ActiveWorkbook

# Rename the 17 "SB Pool N - XXX" sheets to "SB Pool (N) - XXX"
$renames = @(
    @{ Old = "SB Pool 1 - HVAC";      New = "SB Pool (1) - HVAC" },
    @{ Old = "SB Pool 2 - Plumbing";  New = "SB Pool (2) - Plumbing" },
    @{ Old = "SB Pool 3 - ElevM";     New = "SB Pool (3) - ElevM" },
    @{ Old = "SB Pool 4 - ElecM";     New = "SB Pool (4) - ElecM" },
    @{ Old = "SB Pool 5 - Jan";       New = "SB Pool (5) - Jan" },
    @{ Old = "SB Pool 6 - LndScp";    New = "SB Pool (6) - LndScp" },
    @{ Old = "SB Pool 7 - FrSysM";    New = "SB Pool (7) - FrSysM" },
    @{ Old = "SB Pool 8 - FrSupr";    New = "SB Pool (8) - FrSupr" },
    @{ Old = "SB Pool 9 - Roof";      New = "SB Pool (9) - Roof" },
    @{ Old = "SB Pool 10 - BldMgmt";  New = "SB Pool (10) - BldMgmt" },
    @{ Old = "SB Pool 11 - Archt";    New = "SB Pool (11) - Archt" },
    @{ Old = "SB Pool 12 - CommS";    New = "SB Pool (12) - CommS" },
    @{ Old = "SB Pool 13 - ElevIns";  New = "SB Pool (13) - ElevIns" },
    @{ Old = "SB Pool 14 - FacMgmt";  New = "SB Pool (14) - FacMgmt" },
    @{ Old = "SB Pool 15 - Pest";     New = "SB Pool (15) - Pest" },
    @{ Old = "SB Pool 16 - WstMgmt";  New = "SB Pool (16) - WstMgmt" },
    @{ Old = "SB Pool 17 - Cemetary"; New = "SB Pool (17) - Cemetary" }
)

foreach ($r in $renames) {
    $wb.Worksheets($r.Old).Name = $r.New
}

# Move the active tab from "SB Pool (1) - HVAC" to "SB Pool (8) - FrSupr"
$wb.Worksheets("SB Pool (8) - FrSupr").Activate()

# Remove the now-unused custom "Normal 2" cell style
$wb.Styles("Normal 2").Delete()
